$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product_price for row 2 (id 2001)
$ws.Range("E2").Value = 5394

# Resize columns G:J (product_usetaxes..product_subdepartment_id) and K (product_discount)
# ColumnWidth 11.6 maps to the OOXML stored width of 12.5 used elsewhere in the sheet.
$ws.Range($ws.Columns.Item(7), $ws.Columns.Item(10)).ColumnWidth = 11.6
$ws.Columns.Item(11).ColumnWidth = 11.6

# Add the new product row (id 2003 - "Cemento gris x kilo")
$newRow = 4
$ws.Cells.Item($newRow, 1).Value = 2003
$ws.Cells.Item($newRow, 2).Value = 2003
$ws.Cells.Item($newRow, 3).Value = 3
$ws.Cells.Item($newRow, 4).Value = "Cemento gris x kilo"
$ws.Cells.Item($newRow, 5).Value = 176.99
$ws.Cells.Item($newRow, 6).Value = "unidad"
$ws.Cells.Item($newRow, 7).Value = "t"
$ws.Cells.Item($newRow, 8).Value = 13
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 2
$ws.Cells.Item($newRow, 11).Value = 0

# Match the highlighted (yellow fill) + 12pt font styling used by the other data rows
$newRowRange = $ws.Range("A4:K4")
$newRowRange.Interior.Color = 65535
$ws.Range("A4:J4").Font.Size = 12

# Move the active selection down to the newly-added row, like Excel does after entry
$ws.Range("A4:XFD4").Select() | Out-Null
